$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "rect_map": reset the far-field phase/amplitude tables (double
# resonator far field) so every B:G cell in the two 6x6 blocks is 0, except
# the first row of the amplitude table (B2, C2) which stays at 1.
# ---------------------------------------------------------------------------
$map = $wb.Worksheets.Item("rect_map")

# Re-use the existing grey "s=6" cell format (taken from the untouched H2
# cell) for every cell in the two blocks so the style index matches what
# Excel would reuse rather than minting a brand new one.
$map.Range("H2").Copy()
$map.Range("B2:G7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$map.Range("B10:G15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$map.Range("B2:G7").Value = 0
$map.Range("B10:G15").Value = 0

$map.Range("B2").Value = 1
$map.Range("C2").Value = 1

$map.Range("C2").Select()

# ---------------------------------------------------------------------------
# Sheet "rect_array": halve-ish the wavelength variable (B14) from 500nm to
# 100nm.
# ---------------------------------------------------------------------------
$arr = $wb.Worksheets.Item("rect_array")
$arr.Range("B14").Formula = "=100*10^-9"

# Leave the workbook with "rect_array" as the active sheet / B14 selected,
# matching the saved selection state in the workbook.
$arr.Range("B14").Select()
